# Generate Report for Handoff
# Updates the localization-status report to point at the newly generated
# xliff/markdown artifacts (new GUID-based file names) and refreshes the
# handoff/handback timestamps, mirroring a fresh CI run.

$wb = $excel.ActiveWorkbook

$newMd        = "6d63531b-c377-4854-bb56-6b0546cfa32e.md"
$newPath      = "e2e\6d63531b-c377-4854-bb56-6b0546cfa32e.md"
$newZhXlf     = "6d63531b-c377-4854-bb56-6b0546cfa32e.fd6b7621bac04a3b6b2391af348f717d54342167.zh-cn.xlf"
$newDeXlf     = "6d63531b-c377-4854-bb56-6b0546cfa32e.fd6b7621bac04a3b6b2391af348f717d54342167.de-de.xlf"

$genDate      = "2016-09-05 13:11:54"
$zhHandoffDt  = "2016-09-05 13:11:50"

# Hyperlink target addresses are untouched by this edit - only the display
# text (and the underlying cell text) move to the new file name.
$hlAddress = "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/f0e08e975422aa7f4bb0714d2942a25f33629f7f/e2e/36f6d583-7647-4858-971b-d2067cd7c611.md"

# ---------------------------------------------------------------------
# Overview sheet
# ---------------------------------------------------------------------
$wsOverview = $wb.Worksheets.Item("Overview")

$wsOverview.Range("A2").Value = $newMd
$wsOverview.Range("B2").Value = $newPath
$wsOverview.Range("G2").Value = $genDate

$wsOverview.Hyperlinks.Delete()
$wsOverview.Hyperlinks.Add($wsOverview.Range("B2"), $hlAddress, "", "", $newPath)

# ---------------------------------------------------------------------
# zh-cn sheet
# ---------------------------------------------------------------------
$wsZhCn = $wb.Worksheets.Item("zh-cn")

$wsZhCn.Range("A2").Value = $newMd
$wsZhCn.Range("G2").Value = $newZhXlf
$wsZhCn.Range("H2").Value = $zhHandoffDt

$wsZhCn.Hyperlinks.Delete()
$wsZhCn.Hyperlinks.Add($wsZhCn.Range("A2"), $hlAddress, "", "", $newMd)

# ---------------------------------------------------------------------
# de-de sheet
# ---------------------------------------------------------------------
$wsDeDe = $wb.Worksheets.Item("de-de")

$wsDeDe.Range("A2").Value = $newMd
$wsDeDe.Range("G2").Value = $newDeXlf
$wsDeDe.Range("H2").Value = $genDate

$wsDeDe.Hyperlinks.Delete()
$wsDeDe.Hyperlinks.Add($wsDeDe.Range("A2"), $hlAddress, "", "", $newMd)
